# "New function added, Now you can delete Data too."
# Simulates using a new "delete row" feature on the list form to remove
# several people (keeping Amareto and Cleber) and then adding three new
# people (John, Mary, Peter) with the "add" feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete unwanted rows bottom-to-top so row numbers of rows not yet
# processed stay stable:
#   row 9 -> Fabio, 30   (duplicate)
#   row 7 -> Lucas, 32
#   row 6 -> Zezinha, 59
#   row 4 -> Jose, 34
#   row 3 -> Talita, 33
#   row 2 -> Fabio, 30
# This leaves only Amareto (was row 5) and Cleber (was row 8), which end
# up on rows 2 and 3 after the shifts.
$ws.Rows("9:9").Delete()
$ws.Rows("7:7").Delete()
$ws.Rows("6:6").Delete()
$ws.Rows("4:4").Delete()
$ws.Rows("3:3").Delete()
$ws.Rows("2:2").Delete()

# Add the new entries under the remaining two rows.
$ws.Range("B4").Value = "John"
$ws.Range("C4").Value = 52
$ws.Range("B5").Value = "Mary"
$ws.Range("C5").Value = 68
$ws.Range("B6").Value = "Peter"
$ws.Range("C6").Value = 28

# Touch column A so the sheet's recorded extent starts at column A (as in
# the saved workbook), then leave it unformatted again.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Bold = $false

# Final selection left on the sheet after the edits.
[void]$ws.Range("E10").Select()
